$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to remain text so numeric-looking values are not
# auto-converted to numbers (prices use "." as thousands separator too).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.362.85"
$ws.Range("D3").Value = "2.933.12"
$ws.Range("D5").Value = "595.14"
$ws.Range("D6").Value = "143.66"
$ws.Range("D9").Value = "6.94"
$ws.Range("D11").Value = "0.438"
$ws.Range("D12").Value = "0.0000224"
$ws.Range("D13").Value = "33.26"
$ws.Range("D15").Value = "3.418.14"
$ws.Range("D16").Value = "61.373.78"
$ws.Range("D17").Value = "2.932.20"
$ws.Range("D19").Value = "434.26"
$ws.Range("D21").Value = "0.673"
$ws.Range("D22").Value = "7.07"
$ws.Range("D23").Value = "81.41"
$ws.Range("D24").Value = "10.92"
$ws.Range("D25").Value = "2.18"
$ws.Range("D29").Value = "2.60"
$ws.Range("D30").Value = "6.90"
$ws.Range("D31").Value = "26.78"
$ws.Range("D32").Value = "0.108"
$ws.Range("D34").Value = "0.0₃0876"
$ws.Range("D36").Value = "5.62"
$ws.Range("D37").Value = "2.97"
$ws.Range("D40").Value = "8.50"
$ws.Range("D41").Value = "42.06"
$ws.Range("D42").Value = "0.280"
$ws.Range("D43").Value = "0.0344"
$ws.Range("D44").Value = "2.697.31"
$ws.Range("D45").Value = "364.81"
$ws.Range("D46").Value = "133.53"
$ws.Range("D48").Value = "23.58"

# Other column updates (Coin name, Link, Volume percentage)
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -3.70%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("E34").Value = "  +2.73%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("E51").Value = "  +0.14%  "
